# Apply the trading-log update described by the commit:
#   - Trade #364 (MarketMaking) is closed (early_exit) in "All Trades" and
#     in the per-strategy "MarketMaking" sheet.
#   - Four brand-new OPEN trades (#393-#396) are appended to "All Trades"
#     and to their respective per-strategy sheets.
#   - The roll-up "Summary" and "Strategy Status" sheets are refreshed to
#     reflect the new totals.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    # Forces literal text storage so date/time-looking strings (e.g.
    # "2026-02-18") are not auto-coerced into date serials by the
    # COM Value setter.
    param($Sheet, $Row, $Col, $Text)
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

function Set-NumCell {
    param($Sheet, $Row, $Col, $Num)
    $Sheet.Cells.Item($Row, $Col).Value = $Num
}

function Set-PlainCell {
    # Plain text (strategy names, sides, statuses, free-text reasons) -
    # none of these are numeric/date-like, so a direct Value assign is
    # safe and avoids forcing a Text number-format on the cell.
    param($Sheet, $Row, $Col, $Text)
    $Sheet.Cells.Item($Row, $Col).Value = $Text
}

function Clear-Cell {
    param($Sheet, $Row, $Col)
    $Sheet.Cells.Item($Row, $Col).ClearContents()
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
Set-NumCell $summary 3 2 1499.34   # Current Capital
Set-NumCell $summary 4 2 0.45      # Total P&L $
Set-NumCell $summary 6 2 364       # Total Trades
Set-NumCell $summary 7 2 139       # Winning Trades
Set-NumCell $summary 9 2 38.19     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
Set-NumCell $status 6 3 99.09      # Capital
Set-NumCell $status 6 4 221        # Trades
Set-NumCell $status 6 5 -0.72      # P&L $
Set-NumCell $status 6 6 -0.91      # P&L %
Set-NumCell $status 6 7 34.39      # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# --- Close out trade #364 (row 365): MarketMaking, early_exit ---
Set-NumCell  $allTrades 365 7  0.88                 # G Exit Price
Set-PlainCell $allTrades 365 8  "CLOSED"              # H Status
Set-NumCell  $allTrades 365 9  21.6015               # I P&L %
Set-NumCell  $allTrades 365 10 0.16                  # J P&L $
Set-NumCell  $allTrades 365 11 99.09                 # K Capital After
Set-PlainCell $allTrades 365 12 "early_exit"          # L Exit Reason
Set-NumCell  $allTrades 365 13 0.19                  # M Duration (min)

# --- New trade #393: momentum / UP (row 394) ---
Set-NumCell  $allTrades 394 1  393
Set-TextCell $allTrades 394 2  "2026-02-18"
Set-TextCell $allTrades 394 3  "01:54:39"
Set-PlainCell $allTrades 394 4  "momentum"
Set-PlainCell $allTrades 394 5  "UP"
Set-NumCell  $allTrades 394 6  0.723675
Set-PlainCell $allTrades 394 8  "OPEN"
Set-NumCell  $allTrades 394 9  0
Set-NumCell  $allTrades 394 10 0
Set-NumCell  $allTrades 394 11 99.62699700270591
Set-NumCell  $allTrades 394 13 0
Set-NumCell  $allTrades 394 14 0
Set-NumCell  $allTrades 394 15 0
Set-NumCell  $allTrades 394 16 0.9
Set-PlainCell $allTrades 394 17 "Upward momentum: 60.000% over 10 samples"

# --- New trade #394: HighProbConvergence / DOWN (row 395) ---
Set-NumCell  $allTrades 395 1  394
Set-TextCell $allTrades 395 2  "2026-02-18"
Set-TextCell $allTrades 395 3  "01:54:39"
Set-PlainCell $allTrades 395 4  "HighProbConvergence"
Set-PlainCell $allTrades 395 5  "DOWN"
Set-NumCell  $allTrades 395 6  0.29
Set-PlainCell $allTrades 395 8  "OPEN"
Set-NumCell  $allTrades 395 9  0
Set-NumCell  $allTrades 395 10 0
Set-NumCell  $allTrades 395 11 100.1931846556633
Set-NumCell  $allTrades 395 13 0
Set-NumCell  $allTrades 395 14 0
Set-NumCell  $allTrades 395 15 0
Set-NumCell  $allTrades 395 16 0.95
Set-PlainCell $allTrades 395 17 "Mean reversion DOWN: price 55.34% above mean (z=4.36)"

# --- New trade #395: MarketMaking / DOWN (row 396) ---
Set-NumCell  $allTrades 396 1  395
Set-TextCell $allTrades 396 2  "2026-02-18"
Set-TextCell $allTrades 396 3  "01:54:41"
Set-PlainCell $allTrades 396 4  "MarketMaking"
Set-PlainCell $allTrades 396 5  "DOWN"
Set-NumCell  $allTrades 396 6  0.31
Set-PlainCell $allTrades 396 8  "OPEN"
Set-NumCell  $allTrades 396 9  0
Set-NumCell  $allTrades 396 10 0
Set-NumCell  $allTrades 396 11 98.93385807314881
Set-NumCell  $allTrades 396 13 0
Set-NumCell  $allTrades 396 14 0
Set-NumCell  $allTrades 396 15 0
Set-NumCell  $allTrades 396 16 0.6
Set-PlainCell $allTrades 396 17 "Normal spread capture: 238 bps"

# --- New trade #396: EMAArbitrage / UP (row 397) ---
Set-NumCell  $allTrades 397 1  396
Set-TextCell $allTrades 397 2  "2026-02-18"
Set-TextCell $allTrades 397 3  "01:54:41"
Set-PlainCell $allTrades 397 4  "EMAArbitrage"
Set-PlainCell $allTrades 397 5  "UP"
Set-NumCell  $allTrades 397 6  0.71
Set-PlainCell $allTrades 397 8  "OPEN"
Set-NumCell  $allTrades 397 9  0
Set-NumCell  $allTrades 397 10 0
Set-NumCell  $allTrades 397 11 100.430616878256
Set-NumCell  $allTrades 397 13 0
Set-NumCell  $allTrades 397 14 0
Set-NumCell  $allTrades 397 15 0
Set-NumCell  $allTrades 397 16 0.9
Set-PlainCell $allTrades 397 17 "EMA:up, RSI:100.0, ROC:60.00% | 2/3 UP"

# ---------------------------------------------------------------------
# momentum sheet - append trade #393 (row 72)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
Set-NumCell  $momentum 72 1  393
Set-TextCell $momentum 72 2  "2026-02-18"
Set-TextCell $momentum 72 3  "01:54:39"
Set-PlainCell $momentum 72 4  "momentum"
Set-PlainCell $momentum 72 5  "UP"
Set-NumCell  $momentum 72 6  0.723675
Set-PlainCell $momentum 72 8  "OPEN"
Set-NumCell  $momentum 72 9  0
Set-NumCell  $momentum 72 10 0
Set-NumCell  $momentum 72 11 99.62699700270591
Set-NumCell  $momentum 72 12 0
Set-NumCell  $momentum 72 13 0
Set-NumCell  $momentum 72 14 0.9
Set-PlainCell $momentum 72 15 "Upward momentum: 60.000% over 10 samples"
Set-NumCell  $momentum 72 17 0

# ---------------------------------------------------------------------
# HighProbConvergence sheet - append trade #394 (row 36)
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
Set-NumCell  $hpc 36 1  394
Set-TextCell $hpc 36 2  "2026-02-18"
Set-TextCell $hpc 36 3  "01:54:39"
Set-PlainCell $hpc 36 4  "HighProbConvergence"
Set-PlainCell $hpc 36 5  "DOWN"
Set-NumCell  $hpc 36 6  0.29
Set-PlainCell $hpc 36 8  "OPEN"
Set-NumCell  $hpc 36 9  0
Set-NumCell  $hpc 36 10 0
Set-NumCell  $hpc 36 11 100.1931846556633
Set-NumCell  $hpc 36 12 0
Set-NumCell  $hpc 36 13 0
Set-NumCell  $hpc 36 14 0.95
Set-PlainCell $hpc 36 15 "Mean reversion DOWN: price 55.34% above mean (z=4.36)"
Set-NumCell  $hpc 36 17 0

# ---------------------------------------------------------------------
# MarketMaking sheet - close trade #364 (row 222) + append #395 (row 244)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
Set-NumCell  $mm 222 7  0.88                  # G Exit Price
Set-PlainCell $mm 222 8  "CLOSED"              # H Status
Set-NumCell  $mm 222 9  21.6015               # I P&L %
Set-NumCell  $mm 222 10 0.16                  # J P&L $
Set-NumCell  $mm 222 11 99.09                 # K Capital After
Set-PlainCell $mm 222 16 "early_exit"          # P Exit Reason
Set-NumCell  $mm 222 17 0.19                  # Q Duration (min)

Set-NumCell  $mm 244 1  395
Set-TextCell $mm 244 2  "2026-02-18"
Set-TextCell $mm 244 3  "01:54:41"
Set-PlainCell $mm 244 4  "MarketMaking"
Set-PlainCell $mm 244 5  "DOWN"
Set-NumCell  $mm 244 6  0.31
Set-PlainCell $mm 244 8  "OPEN"
Set-NumCell  $mm 244 9  0
Set-NumCell  $mm 244 10 0
Set-NumCell  $mm 244 11 98.93385807314881
Set-NumCell  $mm 244 12 0
Set-NumCell  $mm 244 13 0
Set-NumCell  $mm 244 14 0.6
Set-PlainCell $mm 244 15 "Normal spread capture: 238 bps"
Set-NumCell  $mm 244 17 0

# ---------------------------------------------------------------------
# EMAArbitrage sheet - append trade #396 (row 16)
# ---------------------------------------------------------------------
$ema = $wb.Worksheets.Item("EMAArbitrage")
Set-NumCell  $ema 16 1  396
Set-TextCell $ema 16 2  "2026-02-18"
Set-TextCell $ema 16 3  "01:54:41"
Set-PlainCell $ema 16 4  "EMAArbitrage"
Set-PlainCell $ema 16 5  "UP"
Set-NumCell  $ema 16 6  0.71
Set-PlainCell $ema 16 8  "OPEN"
Set-NumCell  $ema 16 9  0
Set-NumCell  $ema 16 10 0
Set-NumCell  $ema 16 11 100.430616878256
Set-NumCell  $ema 16 12 0
Set-NumCell  $ema 16 13 0
Set-NumCell  $ema 16 14 0.9
Set-PlainCell $ema 16 15 "EMA:up, RSI:100.0, ROC:60.00% | 2/3 UP"
Set-NumCell  $ema 16 17 0

Write-Host "Edit script completed."
